{"js": "// The edit removes the trailing \"blank / page-break / copyright footer\"\n// paragraphs that used to follow the \"LOB1036: Geometria Anal\u00edtica\n// (Requisito fraco)\" requirement line, while leaving the blank paragraph\n// and the page-break paragraph that come right after them (and the final\n// sectPr) untouched.\n//\n// Concretely, right after the \"LOB1036...\" paragraph the body used to have:\n//   1. an empty \"Normal\" paragraph\n//   2. an empty \"Normal\" paragraph with pageBreakBefore\n//   3. the \"\u00a9 2020 . Contact: luizeleno@usp.br. ...\" copyright paragraph\n// All three of those paragraphs must be deleted.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet lobIndex = -1;\nlet copyrightIndex = -1;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text.indexOf(\"LOB1036\") !== -1) {\n    lobIndex = i;\n  }\n  if (text.indexOf(\"Contact: luizeleno@usp.br\") !== -1) {\n    copyrightIndex = i;\n  }\n}\n\nif (lobIndex === -1) {\n  throw new Error('Could not find the \"LOB1036\" requirement paragraph.');\n}\nif (copyrightIndex === -1) {\n  throw new Error('Could not find the copyright/footer paragraph.');\n}\nif (copyrightIndex <= lobIndex) {\n  throw new Error(\"Unexpected paragraph ordering while locating the edit target.\");\n}\n\n// Delete every paragraph strictly after the \"LOB1036...\" one, up to and\n// including the copyright paragraph. Iterate from the end backwards so\n// earlier indices stay valid as items are removed.\nfor (let i = copyrightIndex; i > lobIndex; i--) {\n  paragraphs.items[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# The edit removes the trailing \"blank / page-break / copyright footer\"\n# paragraphs that used to follow the \"LOB1036: Geometria Anal\u00edtica\n# (Requisito fraco)\" requirement line, while leaving the blank paragraph\n# and the page-break paragraph that come right after them (and the final\n# sectPr) untouched.\n#\n# Concretely, right after the \"LOB1036...\" paragraph the body used to have:\n#   1. an empty \"Normal\" paragraph\n#   2. an empty \"Normal\" paragraph with pageBreakBefore\n#   3. the \"\u00a9 2020 . Contact: luizeleno@usp.br. ...\" copyright paragraph\n# All three of those paragraphs must be deleted.\n\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n$lobIndex = -1\n$copyIndex = -1\n\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -like \"*LOB1036*\") {\n        $lobIndex = $i\n    }\n    if ($t -like \"*Contact: luizeleno@usp.br*\") {\n        $copyIndex = $i\n    }\n}\n\nif ($lobIndex -eq -1) {\n    throw \"Could not find the 'LOB1036' requirement paragraph.\"\n}\nif ($copyIndex -eq -1) {\n    throw \"Could not find the copyright/footer paragraph.\"\n}\nif ($copyIndex -le $lobIndex) {\n    throw \"Unexpected paragraph ordering while locating the edit target.\"\n}\n\n# Delete every paragraph strictly after the \"LOB1036...\" one, up to and\n# including the copyright paragraph. Walk backwards so earlier indices\n# remain valid as paragraphs are removed.\nfor ($i = $copyIndex; $i -gt $lobIndex; $i--) {\n    $d.Paragraphs.Item($i).Range.Delete()\n}\n"}
